# Add a new "ssim_dual" results column as the new first column (A),
# shifting the existing ssim_nlm / ssim_gnlm / ssim_bm3d columns to B/C/D.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A; this shifts existing A:C -> B:D and all
# anchored cells (e.g. the stray styled cells in column G/O) along with it.
$ws.Columns("A:A").Insert()

# Header for the new column.
$ws.Range("A1").Value = "ssim_dual"

# Give the new header the same look as the other header cells (bold,
# centered, boxed) so it reads as part of the header row.
$ws.Range("A1").Font.Bold = $true
$ws.Range("A1").Borders.LineStyle = 1
$ws.Range("A1").HorizontalAlignment = -4108
$ws.Range("A1").VerticalAlignment = -4160

# New "ssim_dual" data values for rows 2..51.
$dualValues = @(0.92803188000000003,0.92752071000000003,0.85764662999999997,0.88061412999999999,0.85801879000000003,0.88118377999999997,0.90027073999999996,0.88850467,0.89967065999999996,0.91473526999999999,0.91344555000000005,0.89564876000000004,0.84913817999999996,0.91459407000000004,0.90778106999999997,0.87826716999999999,0.88121590999999999,0.87034800000000001,0.88858946999999999,0.89211187999999997,0.87049370000000004,0.87526499999999996,0.86975391000000002,0.90570198999999996,0.87360167,0.88609245000000003,0.84126372000000005,0.90671550000000001,0.86406099000000003,0.87913390000000002,0.84581052000000001,0.88447237999999995,0.80797922,0.83162398999999998,0.87536744,0.87942914000000005,0.89596763999999995,0.81378550000000005,0.88391713000000005,0.90603402,0.83998474999999995,0.87676414999999996,0.86153341000000006,0.91800106999999997,0.96145497000000002,0.86735382999999999,0.89438706000000001,0.88998641000000001,0.91109214000000005,0.94282699000000003)

for ($i = 0; $i -lt $dualValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 1).Value = $dualValues[$i]
}

# Column A is narrower than the B:D data columns.
$ws.Range("A1:A51").ColumnWidth = 10.14

# Match the recorded selection/cursor position after the edit.
$ws.Range("A2").Select() | Out-Null
